# Update "想去人数" (wish-to-go count) figures for the recurring
# "南宁·熊喵M动漫嘉年华" event, which is listed on both the "展览"
# (Exhibition) sheet and the aggregate "全部类型" (All Types) sheet.
$wb2 = $excel.ActiveWorkbook

$wsExpo = $wb2.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1125
$wsExpo.Range("F4").Value = 2562

$wsAll = $wb2.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1125
$wsAll.Range("F6").Value = 2562
